$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7; this shifts the existing data rows
# (old rows 7-35, one week of prices each) down to rows 8-36, preserving
# all of their values and formatting (matches the diff, where every
# existing row's data moved down by one row and a new, most-recent
# week's data was inserted at row 7).
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new week's data.
$ws.Range("A7").Value = 7
$ws.Range("B7").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C7").Value = "Ñuble"
$ws.Range("D7").Value = 44859
$ws.Range("D7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E7").Value = 16
$ws.Range("F7").Value = 100112037
$ws.Range("G7").Value = "Cebollín"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 160
$ws.Range("K7").Value = 7000
$ws.Range("L7").Value = 7500
$ws.Range("M7").Value = 7250
$ws.Range("N7").Value = "$/docena de atados"
$ws.Range("O7").Value = "Provincia de Diguillín"
$ws.Range("P7").Value = 2417
$ws.Range("Q7").Value = 3
$ws.Range("R7").Value = "Hortaliza"
